$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet 1 (position 1, currently "GNG_TO-...") -> becomes RS_TO sheet
# Shrinks from 5 data rows (A1:B5) to 3 data rows (A1:B3)
# ------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Rows.Item(4).Delete()
$ws1.Rows.Item(4).Delete()
$ws1.Cells.Item(2,2).Value = "eyes open"
$ws1.Cells.Item(3,2).Value = "eyes closed"
$ws1.Name = "RS_TO-16515889550649657"

# ------------------------------------------------------------------
# Sheet 2 (position 2, currently "NB_TO-...") -> becomes GNG_TO sheet
# Shrinks from 10 data rows (A1:B10) to 5 data rows (A1:B5)
# ------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Rows.Item(6).Delete()
$ws2.Rows.Item(6).Delete()
$ws2.Rows.Item(6).Delete()
$ws2.Rows.Item(6).Delete()
$ws2.Rows.Item(6).Delete()
$ws2.Cells.Item(2,2).Value = "go_stims-16515889550707877.csv"
$ws2.Cells.Item(3,2).Value = "GNG_stims-1651588955090667.csv"
$ws2.Cells.Item(4,2).Value = "go_stims-16515889550927384.csv"
$ws2.Cells.Item(5,2).Value = "GNG_stims-16515889551064103.csv"
$ws2.Name = "GNG_TO-16515889551075194"

# ------------------------------------------------------------------
# Sheet 3 (position 3, currently "RS_TO-...") -> becomes vSAT_TO sheet
# Grows from 3 data rows (A1:B3) to 5 data rows (A1:B5)
# ------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("A3").Copy()
$ws3.Range("A4:A5").PasteSpecial(-4122)
$ws3.Cells.Item(2,2).Value = "SAT_stims-1651588955112476.csv"
$ws3.Cells.Item(3,2).Value = "SAT_stims-16515889551379476.csv"
$ws3.Cells.Item(4,1).Value = 2
$ws3.Cells.Item(4,2).Value = "vSAT_stims-16515889551748638.csv"
$ws3.Cells.Item(5,1).Value = 3
$ws3.Cells.Item(5,2).Value = "vSAT_stims-16515889551596172.csv"
$ws3.Name = "vSAT_TO-1651588955190023"

# ------------------------------------------------------------------
# Sheet 4 (position 4, currently "TOL_TO-...") -> stays TOL_TO sheet
# Size unchanged (A1:B7); only the stim filenames change
# ------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Cells.Item(2,2).Value = "MM_stims-1651588955221924.csv"
$ws4.Cells.Item(3,2).Value = "ZM_stims-16515889551963995.csv"
$ws4.Cells.Item(4,2).Value = "MM_stims-16515889552370114.csv"
$ws4.Cells.Item(5,2).Value = "ZM_stims-16515889552229319.csv"
$ws4.Cells.Item(6,2).Value = "MM_stims-165158895526239.csv"
$ws4.Cells.Item(7,2).Value = "ZM_stims-16515889552379646.csv"
$ws4.Name = "TOL_TO-16515889552633896"

# ------------------------------------------------------------------
# Sheet 5 (position 5, currently "vSAT_TO-...") -> becomes NB_TO sheet
# Grows from 5 data rows (A1:B5) to 10 data rows (A1:B10)
# ------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("A5").Copy()
$ws5.Range("A6:A10").PasteSpecial(-4122)
$ws5.Cells.Item(2,2).Value = "OB-16515889558320196.csv"
$ws5.Cells.Item(3,2).Value = "TB-16515889564625697.csv"
$ws5.Cells.Item(4,2).Value = "TB-16515889569751148.csv"
$ws5.Cells.Item(5,2).Value = "ZB-match_1-16515889555728264.csv"
$ws5.Cells.Item(6,1).Value = 4
$ws5.Cells.Item(6,2).Value = "ZB-match_1-16515889557215617.csv"
$ws5.Cells.Item(7,1).Value = 5
$ws5.Cells.Item(7,2).Value = "OB-16515889559706826.csv"
$ws5.Cells.Item(8,1).Value = 6
$ws5.Cells.Item(8,2).Value = "TB-1651588957317661.csv"
$ws5.Cells.Item(9,1).Value = 7
$ws5.Cells.Item(9,2).Value = "ZB-match_0-16515889556321492.csv"
$ws5.Cells.Item(10,1).Value = 8
$ws5.Cells.Item(10,2).Value = "OB-16515889563917174.csv"
$ws5.Name = "NB_TO-16515889573311532"
